# Apply the Jul 10 2023 GitHub Actions cryptos-list refresh to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.198.50'
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").Value = '1.863.64'
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.31%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4708'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.81'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06467'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.75'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07692'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.36%  '
$ws.Range("D13").Value = '1.868.43'
$ws.Range("E13").Value = '  -0.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.53'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.32%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6823'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.075'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '268.49'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").Value = '30.193.26'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.33'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.71%  '
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.0000'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000007519'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.43%  '
$ws.Range("D22").Value = '2.120.49'
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.0000'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.182'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.107'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.316'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.36'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.79'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("E29").Value = '  -3.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.372'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09805'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.510'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.230'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.974'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.04690'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.67%  '
$ws.Range("E36").Value = '  -2.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6850'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.709'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01845'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.737'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.50%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.382'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.92%  '
$ws.Range("E42").Value = '  -2.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8381'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.883'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.76'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4061'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.200'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.931'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '918.42'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.16%  '
$ws.Range("E51").Value = '  -0.54%  '
